# Fill previously-empty indicator cells with "no", per commit message
# "add 'no' for empty values".
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cells = @(
    "P3", "V3", "AB3",
    "L4", "V4",
    "L5", "P5", "V5",
    "L6", "P6",
    "P7", "Q7",
    "L8", "P8", "V8", "AB8",
    "L9",
    "L10", "P10", "V10", "AB10",
    "L11", "P11",
    "L12", "P12", "V12",
    "L13",
    "L14",
    "L15",
    "L16"
)

foreach ($addr in $cells) {
    $ws.Range($addr).Value = "no"
}
